# Parallel Test Report Generation
# Update the Org_Id values produced by the latest parallel test run
# for the "Org" and "Contacts" test script sheets.

$wb = $excel.ActiveWorkbook

$wsOrg = $wb.Worksheets.Item("Org")
$wsOrg.Range("E2").Value = "ACC40211"
$wsOrg.Range("E8").Value = "ACC40209"

$wsContacts = $wb.Worksheets.Item("Contacts")
$wsContacts.Range("E2").Value = "CON23858"
$wsContacts.Range("E5").Value = "CON23859"
$wsContacts.Range("E8").Value = "CON23862"
